# Insert a new column before column A to hold a sequential "id" column.
# The existing user_id / content / group_id columns shift right by one
# (A->B, B->C, C->D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Insert()

# New header
$ws.Range("A1").Value = "id"

# Find the last used data row (originally rows 2..10 held records).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

$ws.Range("D16").Select()
